# Weekly update: insert a new price record as row 41, pushing the
# previously existing rows 41-51 down to 42-52 (the sheet keeps a rolling
# history of weekly entries).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 41; this shifts rows 41:51 down to
# 42:52 and extends the sheet dimension from A1:R51 to A1:R52.
$ws.Rows.Item(41).Insert()

# Populate the new row 41 with this week's data.
$ws.Range("A41").Value = 7
$ws.Range("B41").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C41").Value = "Ñuble"
$ws.Range("D41").Value = 44798
$ws.Range("E41").Value = 16
$ws.Range("F41").Value = 100112013
$ws.Range("G41").Value = "Alcachofa"
$ws.Range("H41").Value = "Madrigal"
$ws.Range("I41").Value = "Primera"
$ws.Range("J41").Value = 50
$ws.Range("K41").Value = 15000
$ws.Range("L41").Value = 15000
$ws.Range("M41").Value = 15000
$ws.Range("N41").Value = "$/caja 40 unidades"
$ws.Range("O41").Value = "Provincia de Limarí"
$ws.Range("P41").Value = 375
$ws.Range("Q41").Value = 40
$ws.Range("R41").Value = "Hortaliza"
